# ajout des différente armes
# Updates weapon names (column J) for several characters' stats rows and
# tweaks a few related numeric stats, then applies the view/print layout
# changes (freeze panes, margins, page setup, header) seen in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New weapon names assigned to the characters' "weaponName" column (J) ---
$ws.Range("J2").Value  = "Hache"           # Stats_Personnage_1 (Djo)
$ws.Range("J4").Value  = "Dague"           # Stats_Personnage_2 (Kral)
$ws.Range("J6").Value  = "BouleElectrique" # Stats_Personnage_3 (Cork)
$ws.Range("J8").Value  = "Epee"            # Stats_Personnage_4 (Nato)
$ws.Range("J10").Value = "Flechette"       # Stats_Personnage_5 (Cyrdin)
$ws.Range("J12").Value = "Lance"           # Stats_Personnage_6 (Galdir)
$ws.Range("J14").Value = "BouleDeFeu"      # Stats_Personnage_7 (Swift)

# --- Related stat tweaks that came along with the new weapons ---
$ws.Range("F6").Value  = 400
$ws.Range("F10").Value = 500
$ws.Range("I10").Value = 20
$ws.Range("F11").Value = 50

# --- Freeze panes on B2 (header row + first column frozen) ---
$ws.Activate()
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Page margins (Normal / narrow metric preset) ---
$ps = $ws.PageSetup
$ps.LeftMargin   = $excel.InchesToPoints(0.70866141732283472)
$ps.RightMargin  = $excel.InchesToPoints(0.70866141732283472)
$ps.TopMargin    = $excel.InchesToPoints(0.74803149606299213)
$ps.BottomMargin = $excel.InchesToPoints(0.74803149606299213)
$ps.HeaderMargin = $excel.InchesToPoints(0.31496062992125984)
$ps.FooterMargin = $excel.InchesToPoints(0.31496062992125984)

# --- Page setup: A4 portrait ---
$ps.PaperSize = 9
$ps.Orientation = 1

# --- Header showing the sheet name ---
$ps.CenterHeader = "&A"
